$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B for the new "StatQuery" column.
# This shifts the existing B (dbExcel/Neo4jData) -> C and C (WebExcel/WebData) -> D.
$ws.Columns("B:B").Insert()

# New header in row 1
$ws.Range("B1").Value = "StatQuery"

# New stat-bar query text in row 2, wrap-text formatted like column A
$ws.Range("B2").Value = 'MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE demo.breed IN[''Samoyed'']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study'
$ws.Range("B2").WrapText = $true

# Match column B's width to column A's width (same visual column sizing)
$ws.Range("B:B").ColumnWidth = $ws.Range("A:A").ColumnWidth

# Update the view: clear the frozen/top-left cell override and move the
# active selection down to A3 (below the two data rows)
$ws.Range("A3").Select()
